# Update CV workbook: add new "SONA Coordinator" position (2025-present)
# for Dr Meek at UoU, below the existing EDUC 6600 Pre-Test Coordinator entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New job entry: row 9 (order / what / when / with / where / why) ---
$ws.Range("A9").Value = 2
$ws.Range("C9").Value = "2025-present"
$ws.Range("B9").Value = "SONA Coordinator"
$ws.Range("D9").Value = "CEHS Office of Research"
$ws.Range("E9").Value = "Utah State University"
$ws.Range("F9").Value = "Administer the participant recruitment system"

# --- Extra bullet points for the new entry (rows 10-11, column F) ---
$ws.Range("F10").Value = "Interface with USU researchers and the IRB"
$ws.Range("F11").Value = "Troubleshoot with USU students and instructors"

# --- Stray formatted cell left from a paste (K16): Segoe UI 12, dark grey ---
$k16Font = $ws.Range("K16").Font
$k16Font.Name = "Segoe UI"
$k16Font.Size = 12
$k16Font.Color = 2696481
$ws.Rows.Item(16).RowHeight = 17.25

# --- Move the active selection to where editing left off ---
$ws.Range("B14").Select() | Out-Null
